$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37/38 swap: RenderToken <-> TrustWalletToken ---
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.04"
$ws.Range("E37").Value = "  -0.83%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'2.35"
$ws.Range("E38").Value = "  +8.25%  "

# --- Row 47/48 swap: WEMIXToken <-> FraxShare ---
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'6.05"
$ws.Range("E47").Value = "  +2.41%  "

$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "'1.08"
$ws.Range("E48").Value = "  +2.85%  "

# --- Price (Column D) updates ---
$ws.Range("D2").Value = "34.045.81"
$ws.Range("D3").Value = "1.787.44"
$ws.Range("D5").Value = "'226.59"
$ws.Range("D8").Value = "'32.15"
$ws.Range("D9").Value = "'0.295"
$ws.Range("D12").Value = "2.046.22"
$ws.Range("D13").Value = "'11.26"
$ws.Range("D14").Value = "1.799.61"
$ws.Range("D15").Value = "34.015.86"
$ws.Range("D17").Value = "'4.18"
$ws.Range("D18").Value = "'67.72"
$ws.Range("D19").Value = "'242.52"
$ws.Range("D22").Value = "'10.70"
$ws.Range("D24").Value = "'2.05"
$ws.Range("D25").Value = "'161.93"
$ws.Range("D27").Value = "'16.20"
$ws.Range("D33").Value = "'3.59"
$ws.Range("D35").Value = "1.397.29"
$ws.Range("D36").Value = "'0.651"
$ws.Range("D40").Value = "'80.17"
$ws.Range("D41").Value = "'2.35"
$ws.Range("D42").Value = "'0.919"
$ws.Range("D43").Value = "'13.71"
$ws.Range("D46").Value = "'0.0509"
$ws.Range("D49").Value = "'107.69"
$ws.Range("D50").Value = "1.946.90"
$ws.Range("D51").Value = "'1.00"

# --- Volume(1h) (Column E) updates ---
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("E6").Value = "  -1.54%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("E9").Value = "  +3.62%  "
$ws.Range("E10").Value = "  -4.45%  "
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("E13").Value = "  +2.84%  "
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("E24").Value = "  -2.69%  "
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("E30").Value = "  +2.28%  "
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("E33").Value = "  +2.51%  "
$ws.Range("E34").Value = "  +1.25%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("E43").Value = "  +13.67%  "
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("E45").Value = "  +8.08%  "
$ws.Range("E46").Value = "  +2.48%  "
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("E51").Value = "  +0.11%  "
